# Updating subpathway analysis section
# The "Number Of Sig Pathways By Model Type" table currently renders every
# run in the Helvetica font family (ascii/hAnsi/eastAsia/cs). Re-point all
# of those runs at Arial instead, without touching anything else (text,
# bold/italic, size, color, borders, etc.).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = ""
$find.Font.Name = "Helvetica"

$find.Replacement.ClearFormatting()
$find.Replacement.Font.Name = "Arial"
$find.Replacement.Font.NameAscii = "Arial"
$find.Replacement.Font.NameFarEast = "Arial"
$find.Replacement.Font.NameBi = "Arial"
$find.Replacement.Font.NameOther = "Arial"

# wdFindContinue (1) / wdReplaceAll (2): hit every run across the document
# (table header + body cells) that currently carries the Helvetica font.
$find.Execute($null, $null, $null, $null, $null, $null, $true, 1, $null, "", 2) | Out-Null
